$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel, so they stay text like the source data.
$textCells = @("D5","D8","D11","D17","D18","D19","D22","D23","D25","D26","D27","D30","D32","D36","D37","D41","D42","D44","D46","D47","D48","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated cell values
$ws.Range("D2").Value = '29.677.68'
$ws.Range("E2").Value = '  +0.44%  '
$ws.Range("D3").Value = '1.606.07'
$ws.Range("E3").Value = '  +0.16%  '
$ws.Range("E4").Value = '  +0.55%  '
$ws.Range("D5").Value = '212.79'
$ws.Range("E5").Value = '  -0.41%  '
$ws.Range("E6").Value = '  +0.74%  '
$ws.Range("E7").Value = '  +0.56%  '
$ws.Range("D8").Value = '28.02'
$ws.Range("E8").Value = '  +3.60%  '
$ws.Range("E9").Value = '  +0.81%  '
$ws.Range("E10").Value = '  +0.65%  '
$ws.Range("D11").Value = '0.0910'
$ws.Range("E11").Value = '  -0.13%  '
$ws.Range("D12").Value = '1.836.16'
$ws.Range("E12").Value = '  +0.19%  '
$ws.Range("D13").Value = '1.611.72'
$ws.Range("E13").Value = '  +0.77%  '
$ws.Range("E14").Value = '  +3.63%  '
$ws.Range("D15").Value = '29.701.89'
$ws.Range("E15").Value = '  +0.36%  '
$ws.Range("E16").Value = '  +0.16%  '
$ws.Range("D17").Value = '64.08'
$ws.Range("E17").Value = '  +0.89%  '
$ws.Range("D18").Value = '241.00'
$ws.Range("E18").Value = '  -0.48%  '
$ws.Range("D19").Value = '7.86'
$ws.Range("E19").Value = '  +3.45%  '
$ws.Range("D20").Value = '0.0₃0697'
$ws.Range("E20").Value = '  +0.62%  '
$ws.Range("E21").Value = '  +0.52%  '
$ws.Range("D22").Value = '4.02'
$ws.Range("E22").Value = '  -0.74%  '
$ws.Range("D23").Value = '9.39'
$ws.Range("E23").Value = '  +1.17%  '
$ws.Range("E24").Value = '  +0.10%  '
$ws.Range("D25").Value = '155.06'
$ws.Range("E25").Value = '  -0.22%  '
$ws.Range("D26").Value = '15.46'
$ws.Range("E26").Value = '  +1.12%  '
$ws.Range("D27").Value = '0.108'
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("E28").Value = '  +0.84%  '
$ws.Range("E29").Value = '  +0.55%  '
$ws.Range("D30").Value = '0.0480'
$ws.Range("E30").Value = '  +1.51%  '
$ws.Range("E31").Value = '  +0.16%  '
$ws.Range("D32").Value = '3.24'
$ws.Range("E32").Value = '  +0.00%  '
$ws.Range("E33").Value = '  +2.27%  '
$ws.Range("D34").Value = '1.428.83'
$ws.Range("E34").Value = '  -0.35%  '
$ws.Range("E35").Value = '  +2.99%  '
$ws.Range("B36").Value = 'MXToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D36").Value = '2.91'
$ws.Range("E36").Value = '  +3.80%  '
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").Value = '1.03'
$ws.Range("E37").Value = '  -1.44%  '
$ws.Range("E38").Value = '  -0.37%  '
$ws.Range("E39").Value = '  +1.89%  '
$ws.Range("E40").Value = '  +2.94%  '
$ws.Range("D41").Value = '56.60'
$ws.Range("E41").Value = '  +5.24%  '
$ws.Range("D42").Value = '0.0498'
$ws.Range("E42").Value = '  +5.97%  '
$ws.Range("E43").Value = '  -0.81%  '
$ws.Range("D44").Value = '0.816'
$ws.Range("E44").Value = '  +1.96%  '
$ws.Range("E45").Value = '  +0.51%  '
$ws.Range("D46").Value = '66.26'
$ws.Range("E46").Value = '  +0.64%  '
$ws.Range("D47").Value = '0.981'
$ws.Range("E47").Value = '  +17.26%  '
$ws.Range("D48").Value = '5.38'
$ws.Range("E48").Value = '  +0.30%  '
$ws.Range("D49").Value = '1.745.64'
$ws.Range("E49").Value = '  +0.13%  '
$ws.Range("D50").Value = '86.68'
$ws.Range("E50").Value = '  +0.23%  '
$ws.Range("E51").Value = '  +2.44%  '
